$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each ticker occupies a contiguous block of 5 rows (A..F).
# Replace ticker names and the B/C/D/E values with the new data.

$data = @(
    # row, ticker, B, C, D, E
    @(2,  "AZO",  1, 0.05971495234788682, 2179.433870967742, 2056.622741935484),
    @(3,  "AZO",  1, 0.1105403963296181,  2420.34935483871,  2179.433870967742),
    @(4,  "AZO",  0, 0.007667381528366324,2438.907096774193, 2420.34935483871),
    @(5,  "AZO",  0, 0.04570489796328658, 2550.377096774193, 2438.907096774193),
    @(6,  "AZO",  0, -0.01588493744568298,2509.864516129032, 2550.377096774193),

    @(7,  "DPZ",  1, 0.02734766921420384, 366.6559349455675, 356.8956702126114),
    @(8,  "DPZ",  0, -0.07257283276264426,340.0466750973318, 366.6559349455675),
    @(9,  "DPZ",  0, -0.04669722367395546,324.1674394507269, 340.0466750973318),
    @(10, "DPZ",  0, -0.04820661875797649,308.5404232833762, 324.1674394507269),
    @(11, "DPZ",  1, 0.2265847966660433,  378.450992356295,  308.5404232833762),

    @(12, "FICO", 1, 0.1545034504137472,  455.323870967742,  394.3893548387097),
    @(13, "FICO", 1, 0.1592521168909191,  527.8351612903226, 455.323870967742),
    @(14, "FICO", 1, 0.2602080072652202,  665.1820967741936, 527.8351612903226),
    @(15, "FICO", 1, 0.1271790964831603,  749.7793548387097, 665.1820967741936),
    @(16, "FICO", 1, 0.1459713309205806,  859.2256451612903, 749.7793548387097),

    @(17, "VRSN", 1, 0.01086814113459522, 185.4240322580645, 183.4304838709677),
    @(18, "VRSN", 1, 0.03492083098112531, 191.8991935483871, 185.4240322580645),
    @(19, "VRSN", 1, 0.07438549305540976, 206.1737096774194, 191.8991935483871),
    @(20, "VRSN", 1, 0.06923460251573021, 220.4480645161291, 206.1737096774194),
    @(21, "VRSN", 0, -0.0592283457884164, 207.3912903225806, 220.4480645161291),

    @(22, "WEC",  1, 0.01539204956007412, 93.1093123853353,  91.69789385851068),
    @(23, "WEC",  0, -0.08907424299135891,84.8156708691656,  93.1093123853353),
    @(24, "WEC",  0, 0.009653156076683933,85.63440977781431, 84.8156708691656),
    @(25, "WEC",  0, 0.006272987585292045,86.17159336722436, 85.63440977781431),
    @(26, "WEC",  0, -0.05177910667074004,81.70970524227522, 86.17159336722436)
)

foreach ($rowData in $data) {
    $r = $rowData[0]
    $ws.Cells.Item($r, 1).Value = $rowData[1]
    $ws.Cells.Item($r, 2).Value = $rowData[2]
    $ws.Cells.Item($r, 3).Value = $rowData[3]
    $ws.Cells.Item($r, 4).Value = $rowData[4]
    $ws.Cells.Item($r, 5).Value = $rowData[5]
}
